$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Oct 29 2020" match row (row 2) is a duplicate/bug entry - remove it so the
# "Oct 7 2020" match (currently row 3) becomes row 2, shrinking the used range
# from A1:K3 to A1:K2.
$ws.Rows.Item(2).Delete()
